$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.056.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.311.02"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.11%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("E7").Value = "  +0.67%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.524"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.76"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.88%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0789"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.76%  "
$ws.Range("E12").Value = "  -1.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.04"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.89"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.669.79"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.284.96"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.789"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.10%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.968.31"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.46"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.60%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0911"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.85%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.18"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.95"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "240.38"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.34%  "
$ws.Range("E24").Value = "  -2.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("E26").Value = "  -1.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.92"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "169.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.06"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -10.31%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.44"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.51%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.22"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.91"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.45%  "
$ws.Range("E34").Value = "  -0.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "18.34"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.94%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.41"
$ws.Range("D36").Style = "Normal"
$ws.Range("E37").Value = "  +0.40%  "
$ws.Range("E38").Value = "  +0.40%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.80"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.88%  "
$ws.Range("E40").Value = "  -2.38%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.110"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.993.80"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.44%  "
$ws.Range("E43").Value = "  +1.04%  "
$ws.Range("E44").Value = "  -0.68%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.11"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -8.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.53"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.84"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.05%  "
$ws.Range("B48").Value = "MultiversX"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.64"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.05%  "
$ws.Range("B49").Value = "BitcoinSV"
$ws.Range("C49").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "75.04"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.535.84"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.54"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.19%  "
